# ============================================================================
# Edit: add "2022-Q3" quarterly sheet + update the "总计" (summary) sheet
# ============================================================================
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Update sheet "总计" (summary): insert a new row for 2022-Q3 at the
#    top of the data block (row 2), shifting all existing rows down by
#    one, and re-numbering the running index in column A.
# ----------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows.Item(2).Insert()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 7
$wsTotal.Range("D2").Value = 1.32
# The newly-inserted row picks up stray formatting from the Insert() -
# clear it on B2:D2 (they should carry the sheet's default style, same
# as every other non-index data cell).
$wsTotal.Range("B2:D2").ClearFormats()

# Re-number the column-A running index (0-based) for every row that
# shifted down.
for ($r = 3; $r -le 9; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}

# Restore the bold/centered/bordered "key" style on A2 (copy it from a
# neighboring column-A cell that still carries it).
$wsTotal.Range("A3").Copy() | Out-Null
$wsTotal.Range("A2").PasteSpecial(-4122) | Out-Null

# ----------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" (i.e.
#    before the existing "2022-Q2" sheet) and populate it.
# ----------------------------------------------------------------------
$wsBefore = $wb.Worksheets.Item("2022-Q2")
$ws = $wb.Worksheets.Add($wsBefore)
$ws.Name = "2022-Q3"

# Match the look & feel of the other quarterly sheets.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# The "fund code" (B) and the numeric-looking measures (D:G) must stay
# TEXT (same as every other quarterly sheet) - set Text format BEFORE
# writing the values so Excel doesn't coerce them (and drop leading
# zeroes on the fund codes).
$ws.Range("B2:B8").NumberFormat = "@"
$ws.Range("D2:G8").NumberFormat = "@"

# Set values row by row for the new "2022-Q3" worksheet
# Row 1
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"
# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "010714"
$ws.Range("C2").Value = "东方红远见价值混合"
$ws.Range("D2").Value = "15.24"
$ws.Range("E2").Value = "94.15"
$ws.Range("F2").Value = "3.56"
$ws.Range("G2").Value = "0.5425"
$ws.Range("H2").Value = 7
# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "013385"
$ws.Range("C3").Value = "信澳优势价值混合A"
$ws.Range("D3").Value = "12.44"
$ws.Range("E3").Value = "84.28"
$ws.Range("F3").Value = "3.77"
$ws.Range("G3").Value = "0.4690"
$ws.Range("H3").Value = 7
# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "000586"
$ws.Range("C4").Value = "景顺长城中小创精选股票"
$ws.Range("D4").Value = "2.21"
$ws.Range("E4").Value = "93.50"
$ws.Range("F4").Value = "5.46"
$ws.Range("G4").Value = "0.1207"
$ws.Range("H4").Value = 10
# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "013393"
$ws.Range("C5").Value = "信澳价值精选混合A"
$ws.Range("D5").Value = "3.34"
$ws.Range("E5").Value = "79.98"
$ws.Range("F5").Value = "3.27"
$ws.Range("G5").Value = "0.1092"
$ws.Range("H5").Value = 9
# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "013386"
$ws.Range("C6").Value = "信澳优势价值混合C"
$ws.Range("D6").Value = "1.26"
$ws.Range("E6").Value = "84.28"
$ws.Range("F6").Value = "3.77"
$ws.Range("G6").Value = "0.0475"
$ws.Range("H6").Value = 7
# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "005444"
$ws.Range("C7").Value = "光大保德信多策略精选18个月定期开放灵活配置混合"
$ws.Range("D7").Value = "0.84"
$ws.Range("E7").Value = "29.28"
$ws.Range("F7").Value = "2.07"
$ws.Range("G7").Value = "0.0174"
$ws.Range("H7").Value = 4
# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "013394"
$ws.Range("C8").Value = "信澳价值精选混合C"
$ws.Range("D8").Value = "0.38"
$ws.Range("E8").Value = "79.98"
$ws.Range("F8").Value = "3.27"
$ws.Range("G8").Value = "0.0124"
$ws.Range("H8").Value = 9

# ----------------------------------------------------------------------
# 3) Apply the bold / centered / bordered "key" style to the header row
#    and to the column-A running index, copying it from the "总计"
#    sheet which already carries it.
# ----------------------------------------------------------------------
$wsTotal.Range("B1").Copy() | Out-Null
$ws.Range("B1:H1").PasteSpecial(-4122) | Out-Null

$wsTotal.Range("A2").Copy() | Out-Null
$ws.Range("A2:A8").PasteSpecial(-4122) | Out-Null
